# WebForm User Assignment execution
# Fills in the PN_Value (phone number) column F for every participant row,
# and records the UnMatchUserPos result (AO2) for row 2.
#
# Values are written as TEXT (not numbers) even though they look numeric,
# matching the original workbook's shared-string storage for this column.
# We briefly force a text number format so Excel stores the literal as a
# string, then restore the cell to the default "Normal" style so no extra
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "F2"  "9840018258"
Set-TextValue "F3"  "9840016718"
Set-TextValue "F4"  "9840092052"
Set-TextValue "F5"  "9840078456"
Set-TextValue "F6"  "9840070290"
Set-TextValue "F7"  "9840026573"
Set-TextValue "F8"  "9840075031"
Set-TextValue "F9"  "9840035059"
Set-TextValue "F10" "9840019442"
Set-TextValue "F11" "9840004019"
Set-TextValue "F12" "9840040903"
Set-TextValue "F13" "9840064219"
Set-TextValue "F14" "9840009933"
Set-TextValue "F15" "9840061523"
Set-TextValue "F16" "9840066966"
Set-TextValue "F17" "9840002096"
Set-TextValue "F18" "9840015635"

Set-TextValue "AO2" "2"
